# Edit script for Cover Letter.docx
# 1) Update the cached DATE field result text (Jan 26 -> Jan 27, 2021)
# 2) Wrap the two "Croptix" occurrences with spellStart/spellEnd proofErr markers
# 3) Split the "Eko-CORE" hyperlink run into "Eko" (wrapped in spellStart/spellEnd) + "-CORE"
# 4) Insert "(acc: 73.68%, AUC: 0.85)" after "...accurately" in the stenosis sentence
#
# Paragraph-scoped Range.InsertXML calls are used so that unrelated runs /
# paragraph marks are left completely untouched; each XML payload is the
# *exact* original paragraph markup (re-extracted byte-for-byte from the
# document) with only the targeted runs modified, so nothing else shifts.

$d = $word.ActiveDocument

# --- 1) the date field result --------------------------------------------
$d.Content.Find.Execute("January 26, 2021", $false, $false, $false, $false, $false, $true, 1, $false, "January 27, 2021", 2) | Out-Null

# --- 2) first "Croptix" occurrence (intro paragraph) ----------------------
$xmlCroptix1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="1BE8724A" w14:textId="7228DC8F" w:rsidR="00927963" w:rsidRPr="00927963" w:rsidRDefault="00927963" w:rsidP="00927963"><w:pPr><w:spacing w:after="0" w:line="331" w:lineRule="atLeast"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">I would like to introduce myself as an applicant for the Data Scientist position at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006133FD" w:rsidRPr="006133FD"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>Croptix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006133FD"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">that was posted on </w:t></w:r><w:r w:rsidR="006133FD"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>Indeed</w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>. As a highly motivated individual who thrives in multidisciplinary teams and is passionate about data science, I believe that my experience and skill set make me the ideal candidate for this position. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(13).Range.InsertXML($xmlCroptix1)

# --- 3) Eko-CORE split + (acc: ...) insertion (same paragraph) -----------
$xmlBody = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="203F1D5C" w14:textId="6B38E09F" w:rsidR="00927963" w:rsidRPr="00927963" w:rsidRDefault="00927963" w:rsidP="00927963"><w:pPr><w:spacing w:after="0" w:line="331" w:lineRule="atLeast"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">During my </w:t></w:r><w:r w:rsidR="00873958"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>internship</w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> at </w:t></w:r><w:r w:rsidR="00873958"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>Ekohealth</w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">, I engaged in extensive data processing, data visualization, and machine learning model development. As a Data Scientist, I demonstrated excellent communication skills, understanding of machine learning/deep learning algorithms, and excellence in Python / R programming.  I also showed proficiency in ggplot2 / </w:t></w:r><w:hyperlink r:id="rId7" w:history="1"><w:r w:rsidR="00284C00" w:rsidRPr="00284C00"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>R shiny</w:t></w:r></w:hyperlink><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>, knowledge of SQL, and excellence in applied statistics, skills that I applied daily while performing exceptionally well.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve">One of </w:t></w:r><w:r w:rsidR="00084C35"><w:t xml:space="preserve">my </w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve">interesting projects has been </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:proofErr w:type="spellStart"/><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Eko</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>-CORE</w:t></w:r></w:hyperlink><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve">, an FDA-cleared digital stethoscope attachment device. I led a research project building the prototype of an audio-based dialysis fistula assessment algorithm to detect stenosis in </w:t></w:r><w:r w:rsidR="00873958"><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve">early stage, which would cost lots of money for patients with arteriovenous fistula (AVF). Using Fast Fourier transform (FFT), my </w:t></w:r><w:r w:rsidR="00873958"><w:t>machine learning</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> model</w:t></w:r><w:r w:rsidR="00873958"><w:t>s</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> detected severe stenosis accurately</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t>(acc: 73.68%, AUC: 0.85)</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> and thus reduced the monthly cost for patients using our products. Ultimately</w:t></w:r><w:r w:rsidR="00873958"><w:t>, the</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> project help</w:t></w:r><w:r w:rsidR="00F26792"><w:t>s</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> secure $295,881 in SBIR research funding for clinical data collection from the National Institutes of Health (NIH).</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="00873958"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve">This experience let me realize how </w:t></w:r><w:r w:rsidR="00873958"><w:t>data</w:t></w:r><w:r w:rsidR="00873958" w:rsidRPr="002119CC"><w:t xml:space="preserve"> science could contribute to people’s quality of life.</w:t></w:r><w:r w:rsidR="00873958"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>I do believe my personality has also played a major role in my successful accomplishments in this field. I am extremely analytical, data-oriented, and calculated. Even in my personal life, I like to analyze the relevant data before making decisions that might optimize outcomes.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(15).Range.InsertXML($xmlBody)

# --- 4) second "Croptix" occurrence (closing paragraph) -------------------
$xmlCroptix2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="3056091F" w14:textId="038721D0" w:rsidR="00927963" w:rsidRPr="00332664" w:rsidRDefault="00927963" w:rsidP="00927963"><w:pPr><w:spacing w:after="0" w:line="331" w:lineRule="atLeast"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>In closing, I believe my education and experience would ensure my success in this Data Scien</w:t></w:r><w:r w:rsidR="00084C35"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>tist</w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> role, and I am enthusiastic to apply these skills within the cutting-edge technological environment I would expect to find at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006133FD" w:rsidRPr="006133FD"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>Croptix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>. Thank you for taking the time to review my application; I look forward to an opportunity to learn more about the Data Scien</w:t></w:r><w:r w:rsidR="00542E8C"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>tist</w:t></w:r><w:r w:rsidRPr="00332664"><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> position and to further discuss my relevant skills and experience. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(17).Range.InsertXML($xmlCroptix2)

Write-Output "done"
